$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out rows 12-15 entirely (table shrinks from A1:F15 to A1:F11)
$ws.Range("A12:F15").Clear()

# Row 6
$ws.Range("E6").Value = "Chemistry_NMC/Graphite"

# Row 7
$ws.Range("C7").Value = "Trigger-Mechanism_Nail"
$ws.Range("E7").Value = "Bottom-Vent-Yes-No"

# Row 8
$ws.Range("A8").Value = "Chemistry_NMC/Graphite"
$ws.Range("C8").Value = "Cell-Capacity-Ah"
$ws.Range("D8").Value = "Chemistry_NMC/Graphite-SiOx"
$ws.Range("E8").Value = "Trigger-Mechanism_Nail"
$ws.Range("F8").Value = "Trigger-Mechanism_Heater (Non-ISC)"

# Row 9
$ws.Range("E9").ClearContents()
$ws.Range("F9").Value = "Chemistry_NMC/Graphite"

# Row 10
$ws.Range("C10").Value = "Trigger-Mechanism_Heater (Non-ISC)"
$ws.Range("E10").ClearContents()
$ws.Range("F10").ClearContents()

# Row 11
$ws.Range("B11").ClearContents()
$ws.Range("C11").Value = "Chemistry_NMC/Graphite"
$ws.Range("E11").ClearContents()
